$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 215-216 (old rows 215-266 shift down to 217-268)
$ws.Rows("215:216").Insert()

# New row 215: Fecha 44754 (2022-07-12), Calidad "Primera"
$ws.Range("A215").Value = 5
$ws.Range("B215").Value = "Macroferia Regional de Talca"
$ws.Range("C215").Value = "Maule"
$ws.Range("D215").Value = 44754
$ws.Range("E215").Value = 7
$ws.Range("F215").Value = 100112008
$ws.Range("G215").Value = "Coliflor"
$ws.Range("H215").Value = "Sin especificar"
$ws.Range("I215").Value = "Primera"
$ws.Range("J215").Value = 2000
$ws.Range("K215").Value = 1000
$ws.Range("L215").Value = 1000
$ws.Range("M215").Value = 1000
$ws.Range("N215").Value = "$/unidad"
$ws.Range("O215").Value = "Región del Maule"
$ws.Range("P215").Value = 1000
$ws.Range("Q215").Value = 1
$ws.Range("R215").Value = "Hortaliza"

# New row 216: Fecha 44754 (2022-07-12), Calidad "Segunda"
$ws.Range("A216").Value = 5
$ws.Range("B216").Value = "Macroferia Regional de Talca"
$ws.Range("C216").Value = "Maule"
$ws.Range("D216").Value = 44754
$ws.Range("E216").Value = 7
$ws.Range("F216").Value = 100112008
$ws.Range("G216").Value = "Coliflor"
$ws.Range("H216").Value = "Sin especificar"
$ws.Range("I216").Value = "Segunda"
$ws.Range("J216").Value = 2000
$ws.Range("K216").Value = 800
$ws.Range("L216").Value = 800
$ws.Range("M216").Value = 800
$ws.Range("N216").Value = "$/unidad"
$ws.Range("O216").Value = "Región del Maule"
$ws.Range("P216").Value = 800
$ws.Range("Q216").Value = 1
$ws.Range("R216").Value = "Hortaliza"
